$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates (new report week / issue number) ---
$ws.Range("A8").Value = "Volume 32   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/27/2025  Through  2/2/2025"

# --- Crime Complaints table: numeric cell updates (rows 14-30) ---
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 7
$ws.Range("H14").Value = -28.571428571428
$ws.Range("I14").Value = 8
$ws.Range("J14").Value = 9
$ws.Range("K14").Value = -11.111111111111
$ws.Range("L14").Value = -27.272727272727
$ws.Range("M14").Value = 14.285714285714
$ws.Range("N14").Value = -82.978723404255
$ws.Range("C15").Value = 12
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = 20
$ws.Range("F15").Value = 31
$ws.Range("G15").Value = 31
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 39
$ws.Range("J15").Value = 36
$ws.Range("K15").Value = 8.333333333333
$ws.Range("L15").Value = 8.333333333333
$ws.Range("M15").Value = 44.444444444444
$ws.Range("N15").Value = -11.363636363636
$ws.Range("C16").Value = 76
$ws.Range("D16").Value = 110
$ws.Range("E16").Value = -30.90909090909
$ws.Range("F16").Value = 308
$ws.Range("H16").Value = -22.025316455696
$ws.Range("I16").Value = 368
$ws.Range("J16").Value = 475
$ws.Range("K16").Value = -22.526315789473
$ws.Range("L16").Value = -11.961722488038
$ws.Range("M16").Value = -3.916449086161
$ws.Range("N16").Value = -76.708860759493
$ws.Range("C17").Value = 146
$ws.Range("D17").Value = 142
$ws.Range("E17").Value = 2.81690140845
$ws.Range("F17").Value = 541
$ws.Range("G17").Value = 529
$ws.Range("H17").Value = 2.26843100189
$ws.Range("I17").Value = 642
$ws.Range("J17").Value = 644
$ws.Range("K17").Value = -0.310559006211
$ws.Range("L17").Value = 4.051863857374
$ws.Range("M17").Value = 87.719298245614
$ws.Range("N17").Value = -6.277372262773
$ws.Range("C18").Value = 51
$ws.Range("E18").Value = -12.068965517241
$ws.Range("F18").Value = 213
$ws.Range("G18").Value = 213
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 246
$ws.Range("J18").Value = 247
$ws.Range("K18").Value = -0.404858299595
$ws.Range("L18").Value = -9.225092250922
$ws.Range("M18").Value = -24.307692307692
$ws.Range("N18").Value = -85.681024447031
$ws.Range("C19").Value = 155
$ws.Range("D19").Value = 222
$ws.Range("E19").Value = -30.18018018018
$ws.Range("F19").Value = 621
$ws.Range("G19").Value = 719
$ws.Range("H19").Value = -13.630041724617
$ws.Range("I19").Value = 697
$ws.Range("J19").Value = 856
$ws.Range("K19").Value = -18.57476635514
$ws.Range("L19").Value = 8.736349453978
$ws.Range("M19").Value = 94.692737430167
$ws.Range("N19").Value = 14.449917898193
$ws.Range("C20").Value = 77
$ws.Range("D20").Value = 83
$ws.Range("E20").Value = -7.22891566265
$ws.Range("F20").Value = 292
$ws.Range("G20").Value = 331
$ws.Range("H20").Value = -11.782477341389
$ws.Range("I20").Value = 350
$ws.Range("J20").Value = 380
$ws.Range("K20").Value = -7.894736842105
$ws.Range("L20").Value = -28.27868852459
$ws.Range("M20").Value = 105.882352941176
$ws.Range("N20").Value = -75.99451303155
$ws.Range("C21").Value = 518
$ws.Range("D21").Value = 626
$ws.Range("E21").Value = -17.252396166134
$ws.Range("F21").Value = 2011
$ws.Range("G21").Value = 2225
$ws.Range("H21").Value = -9.617977528089
$ws.Range("I21").Value = 2350
$ws.Range("J21").Value = 2647
$ws.Range("K21").Value = -11.220249338874
$ws.Range("L21").Value = -5.318291700241
$ws.Range("M21").Value = 45.781637717121
$ws.Range("N21").Value = -61.732616837648
$ws.Range("C22").Value = 7
$ws.Range("D22").Value = 8
$ws.Range("E22").Value = -12.5
$ws.Range("F22").Value = 26
$ws.Range("G22").Value = 30
$ws.Range("H22").Value = -13.333333333333
$ws.Range("I22").Value = 30
$ws.Range("J22").Value = 38
$ws.Range("K22").Value = -21.052631578947
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 20
$ws.Range("C23").Value = 39
$ws.Range("E23").Value = 2.631578947368
$ws.Range("F23").Value = 103
$ws.Range("G23").Value = 127
$ws.Range("H23").Value = -18.897637795275
$ws.Range("I23").Value = 123
$ws.Range("J23").Value = 157
$ws.Range("K23").Value = -21.656050955414
$ws.Range("L23").Value = -24.074074074074
$ws.Range("M23").Value = 33.695652173913
$ws.Range("C24").Value = 323
$ws.Range("D24").Value = 337
$ws.Range("E24").Value = -4.154302670623
$ws.Range("F24").Value = 1296
$ws.Range("G24").Value = 1338
$ws.Range("H24").Value = -3.139013452914
$ws.Range("I24").Value = 1457
$ws.Range("J24").Value = 1529
$ws.Range("K24").Value = -4.708960104643
$ws.Range("L24").Value = -0.951733514615
$ws.Range("M24").Value = 31.261261261261
$ws.Range("C25").Value = 103
$ws.Range("D25").Value = 138
$ws.Range("E25").Value = -25.362318840579
$ws.Range("F25").Value = 447
$ws.Range("G25").Value = 604
$ws.Range("H25").Value = -25.993377483443
$ws.Range("I25").Value = 501
$ws.Range("J25").Value = 686
$ws.Range("K25").Value = -26.967930029154
$ws.Range("L25").Value = -27.913669064748
$ws.Range("C26").Value = 162
$ws.Range("D26").Value = 220
$ws.Range("E26").Value = -26.363636363636
$ws.Range("F26").Value = 705
$ws.Range("G26").Value = 759
$ws.Range("H26").Value = -7.114624505928
$ws.Range("I26").Value = 819
$ws.Range("J26").Value = 892
$ws.Range("K26").Value = -8.183856502242
$ws.Range("L26").Value = 0.61425061425
$ws.Range("M26").Value = 0.367647058823
$ws.Range("C27").Value = 12
$ws.Range("D27").Value = 15
$ws.Range("E27").Value = -20
$ws.Range("F27").Value = 40
$ws.Range("G27").Value = 52
$ws.Range("H27").Value = -23.076923076923
$ws.Range("I27").Value = 50
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = -16.666666666666
$ws.Range("L27").Value = -18.032786885245
$ws.Range("C28").Value = 22
$ws.Range("D28").Value = 16
$ws.Range("E28").Value = 37.5
$ws.Range("F28").Value = 84
$ws.Range("G28").Value = 72
$ws.Range("H28").Value = 16.666666666666
$ws.Range("I28").Value = 98
$ws.Range("J28").Value = 82
$ws.Range("K28").Value = 19.512195121951
$ws.Range("L28").Value = -4.854368932038
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = 25
$ws.Range("F29").Value = 13
$ws.Range("G29").Value = 24
$ws.Range("H29").Value = -45.833333333333
$ws.Range("I29").Value = 22
$ws.Range("J29").Value = 29
$ws.Range("K29").Value = -24.137931034482
$ws.Range("L29").Value = -12
$ws.Range("M29").Value = -33.333333333333
$ws.Range("N29").Value = -79.245283018867
$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 4
$ws.Range("E30").Value = 25
$ws.Range("F30").Value = 12
$ws.Range("H30").Value = -40
$ws.Range("I30").Value = 18
$ws.Range("J30").Value = 24
$ws.Range("K30").Value = -25
$ws.Range("L30").Value = -10
$ws.Range("M30").Value = -40
$ws.Range("N30").Value = -81.632653061224
$ws.Range("G33").Value = 3

# --- Cells that must hold literal text (e.g. "0" or "***.*") ---
# Temporarily force text format so the numeric-looking strings are not
# auto-converted back into numbers, then restore General formatting.
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "***.*"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "***.*"
